# Insert two new "problem" rows at the top of the tracked-problems table
# (row 3), pushing the existing rows down by two, and record two new
# LeetCode entries there: "989. Add to Array-Form of Integer" and
# "415. Add Strings".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert two blank rows before row 3 -----------------
$ws.Rows("3:4").Insert()

# The freshly inserted rows inherit formatting from the row above (no
# border). Copy the cell formatting from row 5 (an untouched data row,
# now shifted down into place) onto the two new rows so they pick up the
# same bordered/wrapped style used by every other entry in the table.
$ws.Range("A5:G5").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122)
$ws.Range("A5:G5").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Fill in the new row 3 : "989. Add to Array-Form of Integer" ----
$ws.Range("A3").Value = "989. Add to Array-Form of Integer"
$ws.Range("B3").Value = "Easy"
$ws.Range("C3").Value = "Array"
$ws.Range("D3").Value = "- Given an array num represent a number and an integer k`n- return array of num + k"
$ws.Range("E3").Value = "- Can use add string algo`n- use the k like a carry`n- loop the num and add to number k, add k % 10 to ans`n- k = / 10 like carry`n- After out the loop, check if k > 0, if yes, add k % 10 to list and k /= 10`n- can use LinkedList to lower add to List to O(1)"
$ws.Rows("3").RowHeight = 112.5

# F3 uses a distinct bold/larger font style (same look as the "Think
# more" callout column elsewhere), so it is left blank but styled.
$ws.Range("F3").Font.Bold = $true
$ws.Range("F3").Font.Size = 14
$ws.Range("F3").HorizontalAlignment = -4131
$ws.Range("F3").VerticalAlignment = -4108
$ws.Range("F3").WrapText = $true

# --- 3. Fill in the new row 4 : "415. Add Strings" ----------------------
$ws.Range("A4").Value = "415. Add Strings"
$ws.Range("B4").Value = "Easy"
$ws.Range("C4").Value = "String"
$ws.Range("D4").Value = "- Given two non-negative integers, num1 and num2 represented as string, return the sum of num1 and num2 as a string."
$ws.Range("E4").Value = "- Think about how to plus two number in primary school`n- take last number of each String`n- cal carry, cal value`n- append value to StringBuilder`n- reverse stringbuilder"
$ws.Range("F4").Value = "Think more"
$ws.Rows("4").RowHeight = 93.75

# --- 4. Fix up conditional formatting ranges (Insert doesn't auto-shift) ---
$cfBlock2 = $ws.Range("A8:G11").FormatConditions
for ($i = 1; $i -le $cfBlock2.Count; $i++) {
    $cfBlock2.Item($i).ModifyAppliesToRange($ws.Range("A10:G13"))
}
$cfBlock2.Item(1).Formula1 = '=$B10="Hard"'
$cfBlock2.Item(2).Formula1 = '=$B10="Medium"'
$cfBlock2.Item(3).Formula1 = '=$B10="Easy"'

$cfBlock4 = $ws.Range("D7").FormatConditions
for ($i = 1; $i -le $cfBlock4.Count; $i++) {
    $cfBlock4.Item($i).ModifyAppliesToRange($ws.Range("D9"))
}
$cfBlock4.Item(1).Formula1 = '=$B9="Hard"'
$cfBlock4.Item(2).Formula1 = '=$B9="Medium"'
$cfBlock4.Item(3).Formula1 = '=$B9="Easy"'

# Blocks 1 and 3 cover multi-area ranges that already include the newly
# inserted rows (A1:G1 + A12:G1048576 automatically grows to cover row 3
# / row 4 too, since A12:G1048576 already started below the insertion
# point... adjust explicitly to match the saved ranges).
$cfBlock1 = $ws.Range("A1:G1").FormatConditions
$cfBlock1.Item(1).ModifyAppliesToRange($ws.Range("A1:G1"))

$cfBlock3 = $ws.Range("A2:G6").FormatConditions
for ($i = 1; $i -le $cfBlock3.Count; $i++) {
    $cfBlock3.Item($i).ModifyAppliesToRange($ws.Range("A2:G8"))
}

# --- 5. Move the active selection to E3, matching the saved view state ---
$ws.Range("E3").Select()
